$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.909.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.78%  "

$ws.Range("D3").Value = "'1.708.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.59%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").Value = "'311.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.44%  "

$ws.Range("D6").Value = "'0.9992"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").Value = "'0.3749"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.11%  "

$ws.Range("D8").Value = "'49.41"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.57%  "

$ws.Range("D9").Value = "'0.3438"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.14%  "

$ws.Range("D10").Value = "'1.215"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.14%  "

$ws.Range("D11").Value = "'0.07535"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.56%  "

$ws.Range("E12").Value = "  +0.08%  "

$ws.Range("D13").Value = "'21.23"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.14%  "

$ws.Range("D14").Value = "'6.304"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.30%  "

$ws.Range("D15").Value = "'7.071"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.66%  "

$ws.Range("D16").Value = "'1.705.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.59%  "

$ws.Range("E17").Value = "  +1.89%  "

$ws.Range("D18").Value = "'0.06723"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.40%  "

$ws.Range("D19").Value = "'0.9987"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("D20").Value = "'84.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.94%  "

$ws.Range("D21").Value = "'17.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.81%  "

$ws.Range("D22").Value = "'6.374"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.19%  "

$ws.Range("D23").Value = "'13.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.58%  "

$ws.Range("D24").Value = "'24.869.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.82%  "

$ws.Range("E25").Value = "  +0.18%  "

$ws.Range("D26").Value = "'2.780"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.29%  "

$ws.Range("D27").Value = "'20.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.70%  "

$ws.Range("D28").Value = "'150.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.75%  "

$ws.Range("D29").Value = "'133.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.81%  "

$ws.Range("D30").Value = "'1.894.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.56%  "

$ws.Range("D31").Value = "'1.230"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +25.83%  "

$ws.Range("D32").Value = "'6.849"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.70%  "

$ws.Range("D33").Value = "'4.241"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.71%  "

$ws.Range("D34").Value = "'13.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.84%  "

$ws.Range("D35").Value = "'0.08796"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.83%  "

$ws.Range("D36").Value = "'1.776"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.88%  "

$ws.Range("D37").Value = "'5.624"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.43%  "

$ws.Range("D38").Value = "'0.06642"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.89%  "

$ws.Range("D39").Value = "'9.160"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.88%  "

$ws.Range("D40").Value = "'0.02407"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.05%  "

$ws.Range("D41").Value = "'0.2224"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.10%  "

$ws.Range("D42").Value = "'1.279"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.01%  "

$ws.Range("E43").Value = "  +4.14%  "

$ws.Range("D44").Value = "'0.9993"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("E45").Value = "  +3.92%  "

$ws.Range("D46").Value = "'0.6145"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.99%  "

$ws.Range("D47").Value = "'3.821"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.26%  "

$ws.Range("D48").Value = "'2.120"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.46%  "

$ws.Range("D49").Value = "'129.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.82%  "

$ws.Range("D50").Value = "'0.07308"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.51%  "

$ws.Range("D51").Value = "'79.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.10%  "
